# Auto-generated edit script: restores FFXIV market-board derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) to match refreshed market data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1185.3396
$ws.Range("I15").Value = 1185.3396
$ws.Range("K15").Value = 3556.0188
$ws.Range("M15").Value = -3387.0188
# row 64
$ws.Range("H64").Value = 46916.523
$ws.Range("I64").Value = 252425
$ws.Range("J64").Value = 3651.5789
$ws.Range("K64").Value = 252425
$ws.Range("L64").Value = 3651.5789
$ws.Range("M64").Value = -252177
$ws.Range("N64").Value = -4147.5789
# row 67
$ws.Range("H67").Value = 46916.523
$ws.Range("I67").Value = 252425
$ws.Range("J67").Value = 3651.5789
$ws.Range("K67").Value = 252425
$ws.Range("L67").Value = 3651.5789
$ws.Range("M67").Value = -251567
$ws.Range("N67").Value = -5367.5789
# row 137
$ws.Range("H137").Value = 1277.119
$ws.Range("I137").Value = 1016.8158
$ws.Range("J137").Value = 3750
$ws.Range("K137").Value = 3050.4474
$ws.Range("L137").Value = 11250
$ws.Range("M137").Value = -500.4474
$ws.Range("N137").Value = -16350

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 2140
$ws.Range("I61").Value = 2014.4286
$ws.Range("K61").Value = 2014.4286
$ws.Range("M61").Value = -1802.4286
# row 74
$ws.Range("H74").Value = 1443.0968
$ws.Range("I74").Value = 782.8889
$ws.Range("J74").Value = 2357.2307
$ws.Range("K74").Value = 782.8889
$ws.Range("L74").Value = 2357.2307
$ws.Range("M74").Value = 91.11109999999996
$ws.Range("N74").Value = -4105.2307
# row 77
$ws.Range("H77").Value = 1443.0968
$ws.Range("I77").Value = 782.8889
$ws.Range("J77").Value = 2357.2307
$ws.Range("K77").Value = 3914.4445
$ws.Range("L77").Value = 11786.1535
$ws.Range("M77").Value = 453.5554999999999
$ws.Range("N77").Value = -20522.1535
# row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
# row 122
$ws.Range("H122").Value = 1525.1875
$ws.Range("I122").Value = 1556.4286
$ws.Range("J122").Value = 1306.5
$ws.Range("K122").Value = 4669.2858
$ws.Range("L122").Value = 3919.5
$ws.Range("M122").Value = -2219.2858
$ws.Range("N122").Value = -8819.5
# row 136
$ws.Range("H136").Value = 2140
$ws.Range("I136").Value = 2014.4286
$ws.Range("K136").Value = 6043.2858
$ws.Range("M136").Value = -3493.2858

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 12789.5
$ws.Range("I58").Value = 2667
$ws.Range("J58").Value = 20381.375
$ws.Range("K58").Value = 2667
$ws.Range("L58").Value = 20381.375
$ws.Range("M58").Value = -2464
$ws.Range("N58").Value = -20787.375
# row 122
$ws.Range("H122").Value = 282.33334
$ws.Range("I122").Value = 282.33334
$ws.Range("K122").Value = 847.0000200000001
$ws.Range("M122").Value = 1602.99998
# row 132
$ws.Range("H132").Value = 3394.0278
$ws.Range("I132").Value = 2976.1155
$ws.Range("J132").Value = 4480.6
$ws.Range("K132").Value = 8928.3465
$ws.Range("L132").Value = 13441.8
$ws.Range("M132").Value = -6398.3465
$ws.Range("N132").Value = -18501.8
# row 134
$ws.Range("H134").Value = 1225.4117
$ws.Range("I134").Value = 1138
$ws.Range("J134").Value = 1633.3334
$ws.Range("K134").Value = 3414
$ws.Range("L134").Value = 4900.0002
$ws.Range("M134").Value = -879
$ws.Range("N134").Value = -9970.0002
# row 136
$ws.Range("H136").Value = 12789.5
$ws.Range("I136").Value = 2667
$ws.Range("J136").Value = 20381.375
$ws.Range("K136").Value = 8001
$ws.Range("L136").Value = 61144.125
$ws.Range("M136").Value = -5451
$ws.Range("N136").Value = -66244.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 69
$ws.Range("H69").Value = 2542.8
$ws.Range("I69").Value = 1500
$ws.Range("J69").Value = 2803.5
$ws.Range("K69").Value = 4500
$ws.Range("L69").Value = 8410.5
$ws.Range("M69").Value = -3689
$ws.Range("N69").Value = -10032.5
# row 72
$ws.Range("H72").Value = 2542.8
$ws.Range("I72").Value = 1500
$ws.Range("J72").Value = 2803.5
$ws.Range("K72").Value = 13500
$ws.Range("L72").Value = 25231.5
$ws.Range("M72").Value = -9444
$ws.Range("N72").Value = -33343.5
# row 131
$ws.Range("H131").Value = 823.87
$ws.Range("J131").Value = 857.90216
$ws.Range("L131").Value = 2573.70648
$ws.Range("N131").Value = -12653.70648

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 93385.87
$ws.Range("I70").Value = 139131.27
$ws.Range("J70").Value = 7613.25
$ws.Range("K70").Value = 139131.27
$ws.Range("L70").Value = 7613.25
$ws.Range("M70").Value = -138861.27
$ws.Range("N70").Value = -8153.25
# row 73
$ws.Range("H73").Value = 93385.87
$ws.Range("I73").Value = 139131.27
$ws.Range("J73").Value = 7613.25
$ws.Range("K73").Value = 139131.27
$ws.Range("L73").Value = 7613.25
$ws.Range("M73").Value = -138195.27
$ws.Range("N73").Value = -9485.25
# row 113
$ws.Range("H113").Value = 2304.7307
$ws.Range("I113").Value = 2665.7273
$ws.Range("J113").Value = 2040
$ws.Range("K113").Value = 2665.7273
$ws.Range("L113").Value = 2040
$ws.Range("M113").Value = -495.7273
$ws.Range("N113").Value = -6380
# row 122
$ws.Range("H122").Value = 1464.871
$ws.Range("I122").Value = 1417.24
$ws.Range("J122").Value = 1663.3334
$ws.Range("K122").Value = 4251.72
$ws.Range("L122").Value = 4990.0002
$ws.Range("M122").Value = -1801.72
$ws.Range("N122").Value = -9890.0002
# row 132
$ws.Range("H132").Value = 2440.6667
$ws.Range("I132").Value = 1873.2727
$ws.Range("J132").Value = 3332.2856
$ws.Range("K132").Value = 5619.8181
$ws.Range("L132").Value = 9996.856800000001
$ws.Range("M132").Value = -3089.8181
$ws.Range("N132").Value = -15056.8568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2788.8696
$ws.Range("I7").Value = 1717.0714
$ws.Range("K7").Value = 1717.0714
$ws.Range("M7").Value = -1605.0714
# row 40
$ws.Range("H40").Value = 93134.45
$ws.Range("I40").Value = 334826.66
$ws.Range("J40").Value = 2499.875
$ws.Range("K40").Value = 334826.66
$ws.Range("L40").Value = 2499.875
$ws.Range("M40").Value = -334690.66
$ws.Range("N40").Value = -2771.875
# row 74
$ws.Range("H74").Value = 17250
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14002
# row 77
$ws.Range("H77").Value = 17250
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40008
# row 122
$ws.Range("H122").Value = 2200
$ws.Range("I122").Value = 2200
$ws.Range("K122").Value = 6600
$ws.Range("M122").Value = -4150
# row 126
$ws.Range("H126").Value = 2788.8696
$ws.Range("I126").Value = 1717.0714
$ws.Range("K126").Value = 5151.2142
$ws.Range("M126").Value = -2681.2142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 1182.3334
$ws.Range("I113").Value = 1099.5
$ws.Range("K113").Value = 3298.5
$ws.Range("M113").Value = -1128.5
# row 133
$ws.Range("H133").Value = 55853.75
$ws.Range("J133").Value = 55853.75
$ws.Range("L133").Value = 55853.75
$ws.Range("N133").Value = -65973.75
